$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = 6
$ws.Range("F7").Value = -15
$ws.Range("F8").Value = -5
